$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06674100000000001
$ws.Range("N2").Value = 0.200223
$ws.Range("O2").Value = 0.001121358778383549
$ws.Range("P2").Value = 0.001121358778383549
$ws.Range("Q2").Value = 13.184242079417
$ws.Range("R2").Value = 118.658178714753
$ws.Range("S2").Value = 0.0003799952710395794
$ws.Range("T2").Value = 0.0003799952710395793
$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("O3").Value = 0.0189158439903152
$ws.Range("P3").Value = 0.01891584399031519
$ws.Range("Q3").Value = 222.4007794046975
$ws.Range("R3").Value = 2001.607014642278
$ws.Range("S3").Value = 0.006410019168355471
$ws.Range("T3").Value = 0.006410019168355469
$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("M4").Value = 31.745291
$ws.Range("N4").Value = 95.235873
$ws.Range("O4").Value = 0.5333731999099544
$ws.Range("P4").Value = 0.5333731999099544
$ws.Range("Q4").Value = 6271.071776352433
$ws.Range("R4").Value = 56439.6459871719
$ws.Range("S4").Value = 0.1807443768864015
$ws.Range("T4").Value = 0.1807443768864014
$ws.Range("G5").Value = 197.5433703333333
$ws.Range("H5").Value = 592.6301109999999
$ws.Range("I5").Value = 0.3388703761585983
$ws.Range("J5").Value = 0.3388703761585982
$ws.Range("M5").Value = 0.07967600000000001
$ws.Range("N5").Value = 0.239028
$ws.Range("O5").Value = 0.001338688093173426
$ws.Range("P5").Value = 0.001338688093173426
$ws.Range("Q5").Value = 15.73946557467867
$ws.Range("R5").Value = 141.655190172108
$ws.Range("S5").Value = 0.0004536417376927156
$ws.Range("T5").Value = 0.0004536417376927155
$ws.Range("G6").Value = 197.5433703333333
$ws.Range("H6").Value = 592.6301109999999
$ws.Range("I6").Value = 0.3388703761585983
$ws.Range("J6").Value = 0.3388703761585982
$ws.Range("M6").Value = 26.500431
$ws.Range("N6").Value = 79.501293
$ws.Range("O6").Value = 0.4452509092281735
$ws.Range("P6").Value = 0.4452509092281735
$ws.Range("Q6").Value = 5234.984455025947
$ws.Range("R6").Value = 47114.86009523352
$ws.Range("S6").Value = 0.1508823430951091
$ws.Range("T6").Value = 0.150882343095109
$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.06674100000000001
$ws.Range("N7").Value = 0.200223
$ws.Range("O7").Value = 0.001121358778383549
$ws.Range("P7").Value = 0.001121358778383549
$ws.Range("Q7").Value = 5.327750670226001
$ws.Range("R7").Value = 47.949756032034
$ws.Range("S7").Value = 0.0001535560442359045
$ws.Range("T7").Value = 0.0001535560442359044
$ws.Range("I8").Value = 0.1369374790620155
$ws.Range("J8").Value = 0.1369374790620154
$ws.Range("O8").Value = 0.0189158439903152
$ws.Range("P8").Value = 0.01891584399031519
$ws.Range("S8").Value = 0.002590287990364138
$ws.Range("T8").Value = 0.002590287990364138
$ws.Range("I9").Value = 0.1369374790620155
$ws.Range("J9").Value = 0.1369374790620154
$ws.Range("M9").Value = 31.745291
$ws.Range("N9").Value = 95.235873
$ws.Range("O9").Value = 0.5333731999099544
$ws.Range("P9").Value = 0.5333731999099544
$ws.Range("Q9").Value = 2534.139365633859
$ws.Range("R9").Value = 22807.25429070473
$ws.Range("S9").Value = 0.07303878139490957
$ws.Range("T9").Value = 0.07303878139490956
$ws.Range("I10").Value = 0.1369374790620155
$ws.Range("J10").Value = 0.1369374790620154
$ws.Range("M10").Value = 0.07967600000000001
$ws.Range("N10").Value = 0.239028
$ws.Range("O10").Value = 0.001338688093173426
$ws.Range("P10").Value = 0.001338688093173426
$ws.Range("Q10").Value = 6.360316183469334
$ws.Range("R10").Value = 57.242845651224
$ws.Range("S10").Value = 0.0001833165727295055
$ws.Range("T10").Value = 0.0001833165727295054
$ws.Range("I11").Value = 0.1369374790620155
$ws.Range("J11").Value = 0.1369374790620154
$ws.Range("M11").Value = 26.500431
$ws.Range("N11").Value = 79.501293
$ws.Range("O11").Value = 0.4452509092281735
$ws.Range("P11").Value = 0.4452509092281735
$ws.Range("Q11").Value = 2115.456601212566
$ws.Range("R11").Value = 19039.10941091309
$ws.Range("S11").Value = 0.06097153705977636
$ws.Range("T11").Value = 0.06097153705977634
$ws.Range("G12").Value = 148.824417
$ws.Range("H12").Value = 446.473251
$ws.Range("I12").Value = 0.2552967790580629
$ws.Range("J12").Value = 0.2552967790580629
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.06674100000000001
$ws.Range("N12").Value = 0.200223
$ws.Range("O12").Value = 0.001121358778383549
$ws.Range("P12").Value = 0.001121358778383549
$ws.Range("Q12").Value = 9.932690414997001
$ws.Range("R12").Value = 89.39421373497301
$ws.Range("S12").Value = 0.0002862792842898042
$ws.Range("T12").Value = 0.0002862792842898041
$ws.Range("G13").Value = 148.824417
$ws.Range("H13").Value = 446.473251
$ws.Range("I13").Value = 0.2552967790580629
$ws.Range("J13").Value = 0.2552967790580629
$ws.Range("O13").Value = 0.0189158439903152
$ws.Range("P13").Value = 0.01891584399031519
$ws.Range("Q13").Value = 167.551390256222
$ws.Range("R13").Value = 1507.962512305998
$ws.Range("S13").Value = 0.004829154043892285
$ws.Range("T13").Value = 0.004829154043892283
$ws.Range("G14").Value = 148.824417
$ws.Range("H14").Value = 446.473251
$ws.Range("I14").Value = 0.2552967790580629
$ws.Range("J14").Value = 0.2552967790580629
$ws.Range("M14").Value = 31.745291
$ws.Range("N14").Value = 95.235873
$ws.Range("O14").Value = 0.5333731999099544
$ws.Range("P14").Value = 0.5333731999099544
$ws.Range("Q14").Value = 4724.474425570347
$ws.Range("R14").Value = 42520.26983013313
$ws.Range("S14").Value = 0.1361684599729036
$ws.Range("T14").Value = 0.1361684599729036
$ws.Range("G15").Value = 148.824417
$ws.Range("H15").Value = 446.473251
$ws.Range("I15").Value = 0.2552967790580629
$ws.Range("J15").Value = 0.2552967790580629
$ws.Range("M15").Value = 0.07967600000000001
$ws.Range("N15").Value = 0.239028
$ws.Range("O15").Value = 0.001338688093173426
$ws.Range("P15").Value = 0.001338688093173426
$ws.Range("Q15").Value = 11.857734248892
$ws.Range("R15").Value = 106.719608240028
$ws.Range("S15").Value = 0.0003417627583505556
$ws.Range("T15").Value = 0.0003417627583505556
$ws.Range("G16").Value = 148.824417
$ws.Range("H16").Value = 446.473251
$ws.Range("I16").Value = 0.2552967790580629
$ws.Range("J16").Value = 0.2552967790580629
$ws.Range("M16").Value = 26.500431
$ws.Range("N16").Value = 79.501293
$ws.Range("O16").Value = 0.4452509092281735
$ws.Range("P16").Value = 0.4452509092281735
$ws.Range("Q16").Value = 3943.911193823727
$ws.Range("R16").Value = 35495.20074441355
$ws.Range("S16").Value = 0.1136711229986266
$ws.Range("T16").Value = 0.1136711229986266
$ws.Range("G17").Value = 35.426853
$ws.Range("H17").Value = 106.280559
$ws.Range("I17").Value = 0.06077202683121193
$ws.Range("J17").Value = 0.06077202683121192
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.06674100000000001
$ws.Range("N17").Value = 0.200223
$ws.Range("O17").Value = 0.001121358778383549
$ws.Range("P17").Value = 0.001121358778383549
$ws.Range("Q17").Value = 2.364423596073
$ws.Range("R17").Value = 21.279812364657
$ws.Range("S17").Value = 0.00006814724576734006
$ws.Range("T17").Value = 0.00006814724576734004
$ws.Range("G18").Value = 35.426853
$ws.Range("H18").Value = 106.280559
$ws.Range("I18").Value = 0.06077202683121193
$ws.Range("J18").Value = 0.06077202683121192
$ws.Range("O18").Value = 0.0189158439903152
$ws.Range("P18").Value = 0.01891584399031519
$ws.Range("Q18").Value = 39.884708384598
$ws.Range("R18").Value = 358.962375461382
$ws.Range("S18").Value = 0.001149554178514454
$ws.Range("T18").Value = 0.001149554178514454
$ws.Range("G19").Value = 35.426853
$ws.Range("H19").Value = 106.280559
$ws.Range("I19").Value = 0.06077202683121193
$ws.Range("J19").Value = 0.06077202683121192
$ws.Range("M19").Value = 31.745291
$ws.Range("N19").Value = 95.235873
$ws.Range("O19").Value = 0.5333731999099544
$ws.Range("P19").Value = 0.5333731999099544
$ws.Range("Q19").Value = 1124.635757699223
$ws.Range("R19").Value = 10121.72181929301
$ws.Range("S19").Value = 0.03241417041597711
$ws.Range("T19").Value = 0.03241417041597711
$ws.Range("G20").Value = 35.426853
$ws.Range("H20").Value = 106.280559
$ws.Range("I20").Value = 0.06077202683121193
$ws.Range("J20").Value = 0.06077202683121192
$ws.Range("M20").Value = 0.07967600000000001
$ws.Range("N20").Value = 0.239028
$ws.Range("O20").Value = 0.001338688093173426
$ws.Range("P20").Value = 0.001338688093173426
$ws.Range("Q20").Value = 2.822669939628001
$ws.Range("R20").Value = 25.404029456652
$ws.Range("S20").Value = 0.00008135478871695939
$ws.Range("T20").Value = 0.00008135478871695936
$ws.Range("G21").Value = 35.426853
$ws.Range("H21").Value = 106.280559
$ws.Range("I21").Value = 0.06077202683121193
$ws.Range("J21").Value = 0.06077202683121192
$ws.Range("M21").Value = 26.500431
$ws.Range("N21").Value = 79.501293
$ws.Range("O21").Value = 0.4452509092281735
$ws.Range("P21").Value = 0.4452509092281735
$ws.Range("Q21").Value = 938.8268734736431
$ws.Range("R21").Value = 8449.441861262787
$ws.Range("S21").Value = 0.02705880020223607
$ws.Range("T21").Value = 0.02705880020223606
$ws.Range("G22").Value = 121.3248153333333
$ws.Range("H22").Value = 363.974446
$ws.Range("I22").Value = 0.2081233388901116
$ws.Range("J22").Value = 0.2081233388901115
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 0.3333333333333333
$ws.Range("M22").Value = 0.06674100000000001
$ws.Range("N22").Value = 0.200223
$ws.Range("O22").Value = 0.001121358778383549
$ws.Range("P22").Value = 0.001121358778383549
$ws.Range("Q22").Value = 8.097339500162001
$ws.Range("R22").Value = 72.876055501458
$ws.Range("S22").Value = 0.0002333809330509208
$ws.Range("T22").Value = 0.0002333809330509208
$ws.Range("G23").Value = 121.3248153333333
$ws.Range("H23").Value = 363.974446
$ws.Range("I23").Value = 0.2081233388901116
$ws.Range("J23").Value = 0.2081233388901115
$ws.Range("O23").Value = 0.0189158439903152
$ws.Range("P23").Value = 0.01891584399031519
$ws.Range("Q23").Value = 136.5914403795676
$ws.Range("R23").Value = 1229.322963416108
$ws.Range("S23").Value = 0.003936828609188849
$ws.Range("T23").Value = 0.003936828609188849
$ws.Range("G24").Value = 121.3248153333333
$ws.Range("H24").Value = 363.974446
$ws.Range("I24").Value = 0.2081233388901116
$ws.Range("J24").Value = 0.2081233388901115
$ws.Range("M24").Value = 31.745291
$ws.Range("N24").Value = 95.235873
$ws.Range("O24").Value = 0.5333731999099544
$ws.Range("P24").Value = 0.5333731999099544
$ws.Range("Q24").Value = 3851.491568277928
$ws.Range("R24").Value = 34663.42411450136
$ws.Range("S24").Value = 0.1110074112397627
$ws.Range("T24").Value = 0.1110074112397626
$ws.Range("G25").Value = 121.3248153333333
$ws.Range("H25").Value = 363.974446
$ws.Range("I25").Value = 0.2081233388901116
$ws.Range("J25").Value = 0.2081233388901115
$ws.Range("M25").Value = 0.07967600000000001
$ws.Range("N25").Value = 0.239028
$ws.Range("O25").Value = 0.001338688093173426
$ws.Range("P25").Value = 0.001338688093173426
$ws.Range("Q25").Value = 9.666675986498667
$ws.Range("R25").Value = 87.00008387848801
$ws.Range("S25").Value = 0.0002786122356836902
$ws.Range("T25").Value = 0.0002786122356836901
$ws.Range("G26").Value = 121.3248153333333
$ws.Range("H26").Value = 363.974446
$ws.Range("I26").Value = 0.2081233388901116
$ws.Range("J26").Value = 0.2081233388901115
$ws.Range("M26").Value = 26.500431
$ws.Range("N26").Value = 79.501293
$ws.Range("O26").Value = 0.4452509092281735
$ws.Range("P26").Value = 0.4452509092281735
$ws.Range("Q26").Value = 3215.159897328742
$ws.Range("R26").Value = 28936.43907595868
$ws.Range("S26").Value = 0.09266710587242545
$ws.Range("T26").Value = 0.09266710587242544
